# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.5 = 5615.34 pesos`n✅ 5615.34 pesos = 1.49 = 1001.15 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("O10").Value = 3737
$wsTasas.Range("N12").Value = 3759.9
$wsTasas.Range("O12").Value = 670.348
